# Update the cryptocurrency price list on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value.
# Only columns B (Coin), C (Link), D (Price), E (Volume(1h)) ever change.
$updates = @{
    2  = @{ D = '74.721.49'; E = '  +1.21%  ' }
    3  = @{ D = '2.812.58'; E = '  +7.45%  ' }
    4  = @{ E = '  -0.03%  ' }
    5  = @{ D = '186.91'; E = '  +0.09%  ' }
    6  = @{ D = '591.81'; E = '  +1.80%  ' }
    7  = @{ D = '0.999'; E = '  +0.00%  ' }
    8  = @{ E = '  +2.85%  ' }
    9  = @{ E = '  -4.98%  ' }
    10 = @{ D = '2.807.08'; E = '  +7.33%  ' }
    11 = @{ E = '  +4.47%  ' }
    12 = @{ E = '  -1.93%  ' }
    13 = @{ D = '4.86'; E = '  +3.88%  ' }
    14 = @{ D = '3.332.51'; E = '  +7.40%  ' }
    15 = @{ D = '74.745.30'; E = '  +1.32%  ' }
    16 = @{ E = '  -1.68%  ' }
    17 = @{ E = '  +1.58%  ' }
    18 = @{ D = '2.816.73'; E = '  +6.57%  ' }
    19 = @{ D = '9.01'; E = '  -0.61%  ' }
    20 = @{ D = '12.23'; E = '  +3.77%  ' }
    21 = @{ D = '376.05'; E = '  +2.83%  ' }
    22 = @{ E = '  -1.80%  ' }
    23 = @{ E = '  -0.25%  ' }
    24 = @{ E = '  +0.01%  ' }
    25 = @{ D = '70.74'; E = '  +1.37%  ' }
    26 = @{ D = '2.944.25'; E = '  +6.73%  ' }
    27 = @{ D = '4.14'; E = '  +0.56%  ' }
    28 = @{ D = '9.66'; E = '  +3.89%  ' }
    29 = @{ E = '  +9.21%  ' }
    30 = @{ E = '  +0.03%  ' }
    31 = @{ B = 'Fetch.AI'; C = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D = '1.39'; E = '  +0.58%  ' }
    32 = @{ B = 'Bittensor'; C = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D = '511.31'; E = '  -1.83%  ' }
    33 = @{ E = '  -0.35%  ' }
    34 = @{ E = '  +2.24%  ' }
    35 = @{ D = '0.998'; E = '  -0.09%  ' }
    36 = @{ D = '164.15'; E = '  +0.95%  ' }
    37 = @{ D = '19.86'; E = '  +4.14%  ' }
    38 = @{ E = '  -0.39%  ' }
    39 = @{ D = '19.34'; E = '  +0.41%  ' }
    40 = @{ D = '182.45'; E = '  +13.14%  ' }
    41 = @{ E = '  -0.05%  ' }
    42 = @{ E = '  +4.34%  ' }
    43 = @{ D = '4.96'; E = '  +1.19%  ' }
    44 = @{ D = '1.66'; E = '  -0.11%  ' }
    45 = @{ E = '  +2.65%  ' }
    46 = @{ D = '39.98'; E = '  +2.78%  ' }
    47 = @{ D = '0.0864'; E = '  +1.27%  ' }
    48 = @{ E = '  -2.76%  ' }
    49 = @{ D = '0.573'; E = '  +9.90%  ' }
    50 = @{ E = '  +2.81%  ' }
    51 = @{ D = '0.633' }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        # Force a text number format so numeric-looking strings (e.g. "186.91")
        # are stored as text rather than being coerced into floating point
        # numbers, then restore the default "Normal" style so no new cell
        # formatting is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
        $cell.Style = "Normal"
    }
}
